$d = $word.ActiveDocument

function Replace-WithTwoRuns($findText, $firstRunText, $secondRunText) {
    $rng = $d.Content
    $found = $rng.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        return
    }

    $target = $d.Range($rng.Start, $rng.End)

    $xml = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:r>
              <w:rPr>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
                <w:lang w:val="en-US"/>
              </w:rPr>
              <w:t>$firstRunText</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
                <w:lang w:val="en-US"/>
              </w:rPr>
              <w:t xml:space="preserve">$secondRunText</w:t>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
"@

    $target.InsertXML($xml)
}

Replace-WithTwoRuns "50 = 25 gold" "Between 40-50" " = 25 gold"
Replace-WithTwoRuns "100 = 80 gold" "Between 90-100" " = 80 gold"
